# Auto-generated Excel COM-interop script applying the BRVM daily data refresh
# Sheet "Recommandations": refresh rows 2-48 (values + re-ranked order); row 48 is new
# Sheet "Top_YTD": refresh B column for rows 2-11 (row 8 unchanged)

$wb = $excel.ActiveWorkbook

$wsReco = $wb.Worksheets.Item("Recommandations")

# Row 2: BRVM - SERVICES PUBLICS
$wsReco.Cells.Item(2, 1).Value = "BRVM - SERVICES PUBLICS"
$wsReco.Cells.Item(2, 2).Value = 0
$wsReco.Cells.Item(2, 3).Value = 8
$wsReco.Cells.Item(2, 4).Value = 3431.33
$wsReco.Cells.Item(2, 5).Value = 112.53
$wsReco.Cells.Item(2, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(2, 7).Value = "➖ Neutre"

# Row 3: AIR LIQUIDE CI
$wsReco.Cells.Item(3, 1).Value = "AIR LIQUIDE CI"
$wsReco.Cells.Item(3, 2).Value = 0
$wsReco.Cells.Item(3, 3).Value = 4
$wsReco.Cells.Item(3, 4).Value = 2795
$wsReco.Cells.Item(3, 5).Value = 700
$wsReco.Cells.Item(3, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(3, 7).Value = "➖ Neutre"

# Row 4: NEI-CEDA CI
$wsReco.Cells.Item(4, 1).Value = "NEI-CEDA CI"
$wsReco.Cells.Item(4, 2).Value = 0
$wsReco.Cells.Item(4, 3).Value = 4
$wsReco.Cells.Item(4, 4).Value = 2755
$wsReco.Cells.Item(4, 5).Value = 695
$wsReco.Cells.Item(4, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(4, 7).Value = "➖ Neutre"

# Row 5: BRVM - AUTRES SECTEURS
$wsReco.Cells.Item(5, 1).Value = "BRVM - AUTRES SECTEURS"
$wsReco.Cells.Item(5, 2).Value = 0
$wsReco.Cells.Item(5, 3).Value = 4
$wsReco.Cells.Item(5, 4).Value = 2452.03
$wsReco.Cells.Item(5, 5).Value = 603.15
$wsReco.Cells.Item(5, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(5, 7).Value = "➖ Neutre"

# Row 6: BRVM - DISTRIBUTION
$wsReco.Cells.Item(6, 1).Value = "BRVM - DISTRIBUTION"
$wsReco.Cells.Item(6, 2).Value = 0
$wsReco.Cells.Item(6, 3).Value = 4
$wsReco.Cells.Item(6, 4).Value = 2222.5
$wsReco.Cells.Item(6, 5).Value = 583.95
$wsReco.Cells.Item(6, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(6, 7).Value = "➖ Neutre"

# Row 7: BRVM - AGRICULTURE
$wsReco.Cells.Item(7, 1).Value = "BRVM - AGRICULTURE"
$wsReco.Cells.Item(7, 2).Value = 0
$wsReco.Cells.Item(7, 3).Value = 4
$wsReco.Cells.Item(7, 4).Value = 1478.23
$wsReco.Cells.Item(7, 5).Value = 369.62
$wsReco.Cells.Item(7, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(7, 7).Value = "➖ Neutre"

# Row 8: BRVM - TRANSPORT
$wsReco.Cells.Item(8, 1).Value = "BRVM - TRANSPORT"
$wsReco.Cells.Item(8, 2).Value = 0
$wsReco.Cells.Item(8, 3).Value = 4
$wsReco.Cells.Item(8, 4).Value = 1429.47
$wsReco.Cells.Item(8, 5).Value = 364.71
$wsReco.Cells.Item(8, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(8, 7).Value = "➖ Neutre"

# Row 9: BRVM - CONSOMMATION DISCRETIONNAIRE
$wsReco.Cells.Item(9, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$wsReco.Cells.Item(9, 2).Value = 0
$wsReco.Cells.Item(9, 3).Value = 4
$wsReco.Cells.Item(9, 4).Value = 824.42
$wsReco.Cells.Item(9, 5).Value = 219.99
$wsReco.Cells.Item(9, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(9, 7).Value = "➖ Neutre"

# Row 10: BRVM - FINANCES
$wsReco.Cells.Item(10, 1).Value = "BRVM - FINANCES"
$wsReco.Cells.Item(10, 2).Value = 0
$wsReco.Cells.Item(10, 3).Value = 4
$wsReco.Cells.Item(10, 4).Value = 570.48
$wsReco.Cells.Item(10, 5).Value = 143.16
$wsReco.Cells.Item(10, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(10, 7).Value = "➖ Neutre"

# Row 11: BRVM-PRESTIGE
$wsReco.Cells.Item(11, 1).Value = "BRVM-PRESTIGE"
$wsReco.Cells.Item(11, 2).Value = 0
$wsReco.Cells.Item(11, 3).Value = 4
$wsReco.Cells.Item(11, 4).Value = 569.41
$wsReco.Cells.Item(11, 5).Value = 143.04
$wsReco.Cells.Item(11, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(11, 7).Value = "➖ Neutre"

# Row 12: BRVM - SERVICES FINANCIERS
$wsReco.Cells.Item(12, 1).Value = "BRVM - SERVICES FINANCIERS"
$wsReco.Cells.Item(12, 2).Value = 0
$wsReco.Cells.Item(12, 3).Value = 4
$wsReco.Cells.Item(12, 4).Value = 560.67
$wsReco.Cells.Item(12, 5).Value = 140.7
$wsReco.Cells.Item(12, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(12, 7).Value = "➖ Neutre"

# Row 13: BRVM - INDUSTRIELS
$wsReco.Cells.Item(13, 1).Value = "BRVM - INDUSTRIELS"
$wsReco.Cells.Item(13, 2).Value = 0
$wsReco.Cells.Item(13, 3).Value = 4
$wsReco.Cells.Item(13, 4).Value = 509.88
$wsReco.Cells.Item(13, 5).Value = 128.59
$wsReco.Cells.Item(13, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(13, 7).Value = "➖ Neutre"

# Row 14: BRVM - ENERGIE
$wsReco.Cells.Item(14, 1).Value = "BRVM - ENERGIE"
$wsReco.Cells.Item(14, 2).Value = 0
$wsReco.Cells.Item(14, 3).Value = 4
$wsReco.Cells.Item(14, 4).Value = 431.78
$wsReco.Cells.Item(14, 5).Value = 107.3
$wsReco.Cells.Item(14, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(14, 7).Value = "➖ Neutre"

# Row 15: BRVM - TELECOMMUNICATIONS
$wsReco.Cells.Item(15, 1).Value = "BRVM - TELECOMMUNICATIONS"
$wsReco.Cells.Item(15, 2).Value = 0
$wsReco.Cells.Item(15, 3).Value = 4
$wsReco.Cells.Item(15, 4).Value = 388
$wsReco.Cells.Item(15, 5).Value = 96.38
$wsReco.Cells.Item(15, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(15, 7).Value = "➖ Neutre"

# Row 16: BRVM - INDUSTRIE                  (**)
$wsReco.Cells.Item(16, 1).Value = "BRVM - INDUSTRIE                  (**)"
$wsReco.Cells.Item(16, 2).Value = 0
$wsReco.Cells.Item(16, 3).Value = 1
$wsReco.Cells.Item(16, 4).Value = 218.47
$wsReco.Cells.Item(16, 5).Value = 218.47
$wsReco.Cells.Item(16, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(16, 7).Value = "➖ Neutre"

# Row 17: BRVM - INDUSTRIE         (**)
$wsReco.Cells.Item(17, 1).Value = "BRVM - INDUSTRIE         (**)"
$wsReco.Cells.Item(17, 2).Value = 0
$wsReco.Cells.Item(17, 3).Value = 1
$wsReco.Cells.Item(17, 4).Value = 214.39
$wsReco.Cells.Item(17, 5).Value = 214.39
$wsReco.Cells.Item(17, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(17, 7).Value = "➖ Neutre"

# Row 18: BRVM-PRINCIPAL                    (**)
$wsReco.Cells.Item(18, 1).Value = "BRVM-PRINCIPAL                    (**)"
$wsReco.Cells.Item(18, 2).Value = 0
$wsReco.Cells.Item(18, 3).Value = 1
$wsReco.Cells.Item(18, 4).Value = 209.74
$wsReco.Cells.Item(18, 5).Value = 209.74
$wsReco.Cells.Item(18, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(18, 7).Value = "➖ Neutre"

# Row 19: BRVM-PRINCIPAL            (**)
$wsReco.Cells.Item(19, 1).Value = "BRVM-PRINCIPAL            (**)"
$wsReco.Cells.Item(19, 2).Value = 0
$wsReco.Cells.Item(19, 3).Value = 1
$wsReco.Cells.Item(19, 4).Value = 205.89
$wsReco.Cells.Item(19, 5).Value = 205.89
$wsReco.Cells.Item(19, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(19, 7).Value = "➖ Neutre"

# Row 20: BRVM - CONSOMMATION DE BASE          (**)
$wsReco.Cells.Item(20, 1).Value = "BRVM - CONSOMMATION DE BASE          (**)"
$wsReco.Cells.Item(20, 2).Value = 0
$wsReco.Cells.Item(20, 3).Value = 1
$wsReco.Cells.Item(20, 4).Value = 193.64
$wsReco.Cells.Item(20, 5).Value = 193.64
$wsReco.Cells.Item(20, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(20, 7).Value = "➖ Neutre"

# Row 21: BRVM - CONSOMMATION DE BASE         (**)
$wsReco.Cells.Item(21, 1).Value = "BRVM - CONSOMMATION DE BASE         (**)"
$wsReco.Cells.Item(21, 2).Value = 0
$wsReco.Cells.Item(21, 3).Value = 1
$wsReco.Cells.Item(21, 4).Value = 190.31
$wsReco.Cells.Item(21, 5).Value = 190.31
$wsReco.Cells.Item(21, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(21, 7).Value = "➖ Neutre"

# Row 22: SAFCA CI (SAFC)
$wsReco.Cells.Item(22, 1).Value = "SAFCA CI (SAFC)"
$wsReco.Cells.Item(22, 2).Value = 4
$wsReco.Cells.Item(22, 3).Value = 0
$wsReco.Cells.Item(22, 4).Value = 28.22
$wsReco.Cells.Item(22, 5).Value = 7.39
$wsReco.Cells.Item(22, 6).Value = "🟢 Achat"
$wsReco.Cells.Item(22, 7).Value = "✅ Renforcer"

# Row 23: CFAO MOTORS CI (CFAC)
$wsReco.Cells.Item(23, 1).Value = "CFAO MOTORS CI (CFAC)"
$wsReco.Cells.Item(23, 2).Value = 3
$wsReco.Cells.Item(23, 3).Value = 0
$wsReco.Cells.Item(23, 4).Value = 22.31
$wsReco.Cells.Item(23, 5).Value = 7.43
$wsReco.Cells.Item(23, 6).Value = "🟢 Achat"
$wsReco.Cells.Item(23, 7).Value = "✅ Renforcer"

# Row 24: SERVAIR ABIDJAN CI (ABJC)
$wsReco.Cells.Item(24, 1).Value = "SERVAIR ABIDJAN CI (ABJC)"
$wsReco.Cells.Item(24, 2).Value = 2
$wsReco.Cells.Item(24, 3).Value = 0
$wsReco.Cells.Item(24, 4).Value = 14.82
$wsReco.Cells.Item(24, 5).Value = 7.5
$wsReco.Cells.Item(24, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(24, 7).Value = "➖ Neutre"

# Row 25: BERNABE CI (BNBC)
$wsReco.Cells.Item(25, 1).Value = "BERNABE CI (BNBC)"
$wsReco.Cells.Item(25, 2).Value = 2
$wsReco.Cells.Item(25, 3).Value = 0
$wsReco.Cells.Item(25, 4).Value = 9.91
$wsReco.Cells.Item(25, 5).Value = 4.4
$wsReco.Cells.Item(25, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(25, 7).Value = "➖ Neutre"

# Row 26: UNILEVER CI (UNLC)
$wsReco.Cells.Item(26, 1).Value = "UNILEVER CI (UNLC)"
$wsReco.Cells.Item(26, 2).Value = 1
$wsReco.Cells.Item(26, 3).Value = 0
$wsReco.Cells.Item(26, 4).Value = 7.5
$wsReco.Cells.Item(26, 5).Value = 7.5
$wsReco.Cells.Item(26, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(26, 7).Value = "➖ Neutre"

# Row 27: SETAO CI (STAC)
$wsReco.Cells.Item(27, 1).Value = "SETAO CI (STAC)"
$wsReco.Cells.Item(27, 2).Value = 1
$wsReco.Cells.Item(27, 3).Value = 0
$wsReco.Cells.Item(27, 4).Value = 7.5
$wsReco.Cells.Item(27, 5).Value = 7.5
$wsReco.Cells.Item(27, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(27, 7).Value = "➖ Neutre"

# Row 28: SUCRIVOIRE (SCRC)
$wsReco.Cells.Item(28, 1).Value = "SUCRIVOIRE (SCRC)"
$wsReco.Cells.Item(28, 2).Value = 1
$wsReco.Cells.Item(28, 3).Value = 0
$wsReco.Cells.Item(28, 4).Value = 5.68
$wsReco.Cells.Item(28, 5).Value = 5.68
$wsReco.Cells.Item(28, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(28, 7).Value = "➖ Neutre"

# Row 29: ECOBANK COTE D''IVOIRE (ECOC)
$wsReco.Cells.Item(29, 1).Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$wsReco.Cells.Item(29, 2).Value = 1
$wsReco.Cells.Item(29, 3).Value = 0
$wsReco.Cells.Item(29, 4).Value = 5.6
$wsReco.Cells.Item(29, 5).Value = 5.6
$wsReco.Cells.Item(29, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(29, 7).Value = "➖ Neutre"

# Row 30: ORAGROUP TOGO (ORGT)
$wsReco.Cells.Item(30, 1).Value = "ORAGROUP TOGO (ORGT)"
$wsReco.Cells.Item(30, 2).Value = 1
$wsReco.Cells.Item(30, 3).Value = 0
$wsReco.Cells.Item(30, 4).Value = 4.55
$wsReco.Cells.Item(30, 5).Value = 4.55
$wsReco.Cells.Item(30, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(30, 7).Value = "➖ Neutre"

# Row 31: NEI-CEDA CI (NEIC)
$wsReco.Cells.Item(31, 1).Value = "NEI-CEDA CI (NEIC)"
$wsReco.Cells.Item(31, 2).Value = 1
$wsReco.Cells.Item(31, 3).Value = 0
$wsReco.Cells.Item(31, 4).Value = 2.86
$wsReco.Cells.Item(31, 5).Value = 2.86
$wsReco.Cells.Item(31, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(31, 7).Value = "➖ Neutre"

# Row 32: ECOBANK TRANS. INCORP. TG (ETIT)
$wsReco.Cells.Item(32, 1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$wsReco.Cells.Item(32, 2).Value = 1
$wsReco.Cells.Item(32, 3).Value = 1
$wsReco.Cells.Item(32, 4).Value = 0.21
$wsReco.Cells.Item(32, 5).Value = -4.55
$wsReco.Cells.Item(32, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(32, 7).Value = "👀 À surveiller"

# Row 33: AFRICA GLOBAL LOGISTICS CI (SDSC)
$wsReco.Cells.Item(33, 1).Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$wsReco.Cells.Item(33, 2).Value = 1
$wsReco.Cells.Item(33, 3).Value = 1
$wsReco.Cells.Item(33, 4).Value = 0.08
$wsReco.Cells.Item(33, 5).Value = 2.76
$wsReco.Cells.Item(33, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(33, 7).Value = "👀 À surveiller"

# Row 34: TOTAL
$wsReco.Cells.Item(34, 1).Value = "TOTAL"
$wsReco.Cells.Item(34, 2).Value = 0
$wsReco.Cells.Item(34, 3).Value = 4
$wsReco.Cells.Item(34, 4).Value = 0
$wsReco.Cells.Item(34, 5).Value = 0
$wsReco.Cells.Item(34, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(34, 7).Value = "➖ Neutre"

# Row 35: TOTALENERGIES MARKETING SN (TTLS)
$wsReco.Cells.Item(35, 1).Value = "TOTALENERGIES MARKETING SN (TTLS)"
$wsReco.Cells.Item(35, 2).Value = 0
$wsReco.Cells.Item(35, 3).Value = 1
$wsReco.Cells.Item(35, 4).Value = -0.99
$wsReco.Cells.Item(35, 5).Value = -0.99
$wsReco.Cells.Item(35, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(35, 7).Value = "➖ Neutre"

# Row 36: AIR LIQUIDE CI (SIVC)
$wsReco.Cells.Item(36, 1).Value = "AIR LIQUIDE CI (SIVC)"
$wsReco.Cells.Item(36, 2).Value = 1
$wsReco.Cells.Item(36, 3).Value = 1
$wsReco.Cells.Item(36, 4).Value = -1.22
$wsReco.Cells.Item(36, 5).Value = -5.63
$wsReco.Cells.Item(36, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(36, 7).Value = "👀 À surveiller"

# Row 37: BANK OF AFRICA ML (BOAM)
$wsReco.Cells.Item(37, 1).Value = "BANK OF AFRICA ML (BOAM)"
$wsReco.Cells.Item(37, 2).Value = 0
$wsReco.Cells.Item(37, 3).Value = 1
$wsReco.Cells.Item(37, 4).Value = -1.37
$wsReco.Cells.Item(37, 5).Value = -1.37
$wsReco.Cells.Item(37, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(37, 7).Value = "➖ Neutre"

# Row 38: TOTALENERGIES MARKETING CI (TTLC)
$wsReco.Cells.Item(38, 1).Value = "TOTALENERGIES MARKETING CI (TTLC)"
$wsReco.Cells.Item(38, 2).Value = 0
$wsReco.Cells.Item(38, 3).Value = 1
$wsReco.Cells.Item(38, 4).Value = -1.67
$wsReco.Cells.Item(38, 5).Value = -1.67
$wsReco.Cells.Item(38, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(38, 7).Value = "➖ Neutre"

# Row 39: SMB CI (SMBC)
$wsReco.Cells.Item(39, 1).Value = "SMB CI (SMBC)"
$wsReco.Cells.Item(39, 2).Value = 0
$wsReco.Cells.Item(39, 3).Value = 1
$wsReco.Cells.Item(39, 4).Value = -1.89
$wsReco.Cells.Item(39, 5).Value = -1.89
$wsReco.Cells.Item(39, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(39, 7).Value = "➖ Neutre"

# Row 40: SITAB CI (STBC)
$wsReco.Cells.Item(40, 1).Value = "SITAB CI (STBC)"
$wsReco.Cells.Item(40, 2).Value = 0
$wsReco.Cells.Item(40, 3).Value = 1
$wsReco.Cells.Item(40, 4).Value = -2.02
$wsReco.Cells.Item(40, 5).Value = -2.02
$wsReco.Cells.Item(40, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(40, 7).Value = "➖ Neutre"

# Row 41: BANK OF AFRICA BF (BOABF)
$wsReco.Cells.Item(41, 1).Value = "BANK OF AFRICA BF (BOABF)"
$wsReco.Cells.Item(41, 2).Value = 0
$wsReco.Cells.Item(41, 3).Value = 1
$wsReco.Cells.Item(41, 4).Value = -3.31
$wsReco.Cells.Item(41, 5).Value = -3.31
$wsReco.Cells.Item(41, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(41, 7).Value = "➖ Neutre"

# Row 42: SICABLE CI (CABC)
$wsReco.Cells.Item(42, 1).Value = "SICABLE CI (CABC)"
$wsReco.Cells.Item(42, 2).Value = 0
$wsReco.Cells.Item(42, 3).Value = 1
$wsReco.Cells.Item(42, 4).Value = -4.19
$wsReco.Cells.Item(42, 5).Value = -4.19
$wsReco.Cells.Item(42, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(42, 7).Value = "➖ Neutre"

# Row 43: BANK OF AFRICA SENEGAL (BOAS)
$wsReco.Cells.Item(43, 1).Value = "BANK OF AFRICA SENEGAL (BOAS)"
$wsReco.Cells.Item(43, 2).Value = 0
$wsReco.Cells.Item(43, 3).Value = 2
$wsReco.Cells.Item(43, 4).Value = -4.24
$wsReco.Cells.Item(43, 5).Value = -2.1
$wsReco.Cells.Item(43, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(43, 7).Value = "➖ Neutre"

# Row 44: FILTISAC CI (FTSC)
$wsReco.Cells.Item(44, 1).Value = "FILTISAC CI (FTSC)"
$wsReco.Cells.Item(44, 2).Value = 0
$wsReco.Cells.Item(44, 3).Value = 2
$wsReco.Cells.Item(44, 4).Value = -4.32
$wsReco.Cells.Item(44, 5).Value = -3.12
$wsReco.Cells.Item(44, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(44, 7).Value = "➖ Neutre"

# Row 45: VIVO ENERGY CI (SHEC)
$wsReco.Cells.Item(45, 1).Value = "VIVO ENERGY CI (SHEC)"
$wsReco.Cells.Item(45, 2).Value = 0
$wsReco.Cells.Item(45, 3).Value = 1
$wsReco.Cells.Item(45, 4).Value = -4.8
$wsReco.Cells.Item(45, 5).Value = -4.8
$wsReco.Cells.Item(45, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(45, 7).Value = "➖ Neutre"

# Row 46: LOTERIE NATIONALE DU BENIN (LNBB)
$wsReco.Cells.Item(46, 1).Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$wsReco.Cells.Item(46, 2).Value = 0
$wsReco.Cells.Item(46, 3).Value = 2
$wsReco.Cells.Item(46, 4).Value = -5.05
$wsReco.Cells.Item(46, 5).Value = -1.55
$wsReco.Cells.Item(46, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(46, 7).Value = "➖ Neutre"

# Row 47: UNIWAX CI (UNXC)
$wsReco.Cells.Item(47, 1).Value = "UNIWAX CI (UNXC)"
$wsReco.Cells.Item(47, 2).Value = 0
$wsReco.Cells.Item(47, 3).Value = 1
$wsReco.Cells.Item(47, 4).Value = -5.08
$wsReco.Cells.Item(47, 5).Value = -5.08
$wsReco.Cells.Item(47, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(47, 7).Value = "➖ Neutre"

# Row 48: ORANGE COTE D'IVOIRE (ORAC)
$wsReco.Cells.Item(48, 1).Value = "ORANGE COTE D'IVOIRE (ORAC)"
$wsReco.Cells.Item(48, 2).Value = 0
$wsReco.Cells.Item(48, 3).Value = 2
$wsReco.Cells.Item(48, 4).Value = -5.16
$wsReco.Cells.Item(48, 5).Value = -2.03
$wsReco.Cells.Item(48, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(48, 7).Value = "➖ Neutre"

$wsTop = $wb.Worksheets.Item("Top_YTD")

# Row 2: BRVM - SERVICES PUBLICS
$wsTop.Cells.Item(2, 1).Value = "BRVM - SERVICES PUBLICS"
$wsTop.Cells.Item(2, 2).Value = 10461277.26

# Row 3: AIR LIQUIDE CI
$wsTop.Cells.Item(3, 1).Value = "AIR LIQUIDE CI"
$wsTop.Cells.Item(3, 2).Value = 406682

# Row 4: NEI-CEDA CI
$wsTop.Cells.Item(4, 1).Value = "NEI-CEDA CI"
$wsTop.Cells.Item(4, 2).Value = 386920.11

# Row 5: BRVM - AUTRES SECTEURS
$wsTop.Cells.Item(5, 1).Value = "BRVM - AUTRES SECTEURS"
$wsTop.Cells.Item(5, 2).Value = 258301.54

# Row 6: BRVM - DISTRIBUTION
$wsTop.Cells.Item(6, 1).Value = "BRVM - DISTRIBUTION"
$wsTop.Cells.Item(6, 2).Value = 184294.26

# Row 7: BRVM - AGRICULTURE
$wsTop.Cells.Item(7, 1).Value = "BRVM - AGRICULTURE"
$wsTop.Cells.Item(7, 2).Value = 48513.08

# Row 8: BRVM - TRANSPORT
$wsTop.Cells.Item(8, 1).Value = "BRVM - TRANSPORT"
$wsTop.Cells.Item(8, 2).Value = 43650.91

# Row 9: BRVM - CONSOMMATION DISCRETIONNAIRE
$wsTop.Cells.Item(9, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$wsTop.Cells.Item(9, 2).Value = 8660.66

# Row 10: BRVM - FINANCES
$wsTop.Cells.Item(10, 1).Value = "BRVM - FINANCES"
$wsTop.Cells.Item(10, 2).Value = 3364.94

# Row 11: BRVM-PRESTIGE
$wsTop.Cells.Item(11, 1).Value = "BRVM-PRESTIGE"
$wsTop.Cells.Item(11, 2).Value = 3349.65

